$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells that change, to preserve exact string formatting
$dCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D14","D15","D17","D18","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply cell value updates as described by the diff
$ws.Range('D2').Value = '27.256.97'
$ws.Range('D3').Value = '1.909.49'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '307.72'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '0.5250'
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('D9').Value = '0.07310'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').Value = '21.61'
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('D11').Value = '0.9061'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '0.08236'
$ws.Range('E12').Value = '  -3.36%  '
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').Value = '5.371'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').Value = '1.745.71'
$ws.Range('E15').Value = '  -8.33%  '
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '0.000008683'
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').Value = '14.76'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').Value = '27.291.65'
$ws.Range('E20').Value = '  +0.39%  '
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').Value = '10.82'
$ws.Range('E22').Value = '  +2.05%  '
$ws.Range('D23').Value = '6.500'
$ws.Range('E23').Value = '  +1.19%  '
$ws.Range('D24').Value = '2.341'
$ws.Range('E24').Value = '  +2.31%  '
$ws.Range('D25').Value = '150.01'
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('D26').Value = '18.25'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('D27').Value = '1.735'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = '117.03'
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').Value = '4.855'
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('D30').Value = '4.880'
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').Value = '0.09234'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('E32').Value = '  +2.42%  '
$ws.Range('D33').Value = '0.05080'
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').Value = '1.232'
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D35').Value = '2.990'
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('D36').Value = '3.367'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('D37').Value = '2.741'
$ws.Range('E37').Value = '  +4.79%  '
$ws.Range('D38').Value = '0.5755'
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('D39').Value = '0.02003'
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').Value = '1.084'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').Value = '9.062'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').Value = '6.618'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '116.89'
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('E44').Value = '  +0.34%  '
$ws.Range('D45').Value = '0.4936'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.21'
$ws.Range('E46').Value = '  +0.62%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').Value = '1.644'
$ws.Range('D49').Value = '38.67'
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').Value = '64.06'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').Value = '0.06054'
$ws.Range('E51').Value = '  +2.08%  '
